$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: delete the three paragraphs describing the construction of the
# temporal block matrix (the "The observations are sorted..." paragraph, the
# TT= summation equation paragraph, and the blank paragraph after it) that
# sit right before "Spatial Weight Matrix (".
# ---------------------------------------------------------------------------
$pStart = $null
$pEnd = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("The observations are sorted by closing date")) {
        $pStart = $p
    }
    if ($pStart -ne $null -and $t.StartsWith("Spatial Weight Matrix")) {
        $pEnd = $d.Paragraphs.Item($i - 1)
        break
    }
}
if ($pStart -ne $null -and $pEnd -ne $null) {
    $rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $rng.Delete()
}

# ---------------------------------------------------------------------------
# Change 2: remove the page-break rendering artifact that sits between
# "...sold in the future, " and "but the neighbors...", merging the two runs
# back into one (the text itself is unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "a neighbor to a unit sold in the future, but the neighbors for these units in the seed are not determined.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a neighbor to a unit sold in the future, but the neighbors for these units in the seed are not determined.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: split the run containing "...sales two-days prior, the first
# block becomes a " into two runs at "the " / "first block becomes a ".
# ---------------------------------------------------------------------------
$rng3 = $d.Content.Duplicate
$rng3.Find.Execute("first block becomes a ") | Out-Null
if ($rng3.Find.Found) {
    $tailRange = $d.Range($rng3.Start, $rng3.End)
    $tailRange.Font.Bold = 1
    $tailRange.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 4: rewrite the closing sentence - fix "to includes" -> "to include"
# and change the window from "sixty days" to "120 days, or about four
# months". The author made this edit as several separate text insertions, so
# reproduce the same run boundaries (all runs share identical formatting).
# ---------------------------------------------------------------------------
$fullSentence = "To create the full temporal matrix, we can then simply sum the daily matrix to include all sales that occur for a chosen window of time which, in this paper, is assumed to be 120 days, or about four months."
$d.Content.Find.Execute(
    "To create the full temporal matrix, we can then simply sum the daily matrix to includes all sales that occur for a chosen window of time which, in this paper, is assumed to be sixty days.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $fullSentence,
    2) | Out-Null

$rngFull = $d.Content.Duplicate
$rngFull.Find.Execute($fullSentence) | Out-Null
if ($rngFull.Find.Found) {
    $pieces = @(
        "To create the full temporal matrix, we can then simply sum the daily matrix to ",
        "include",
        " all sales that occur for a chosen window of time which, in this paper, is assumed to be ",
        "120 days, or about four months",
        "."
    )
    $pos = $rngFull.Start
    $bounds = @()
    foreach ($piece in $pieces) {
        $bounds += , @($pos, $pos + $piece.Length)
        $pos += $piece.Length
    }
    foreach ($b in $bounds) {
        $pieceRange = $d.Range($b[0], $b[1])
        $pieceRange.Font.Bold = 1
        $pieceRange.Font.Bold = 0
    }
}
